$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (row, date-serial, B, C, D) appended after the existing last row (357)
# as part of the "aggiornamento fino a 1/09/2021" update.
$rows = @(
    @(358, 44432, 0, 2,  32.34675723758694),
    @(359, 44433, 0, 2,  32.34675723758694),
    @(360, 44434, 0, 2,  32.34675723758694),
    @(361, 44435, 1, 3,  48.5201358563804),
    @(362, 44436, 2, 4,  64.69351447517387),
    @(363, 44437, 2, 5,  80.86689309396733),
    @(364, 44438, 6, 11, 177.9071648067281),
    @(365, 44439, 0, 11, 177.9071648067281),
    @(366, 44440, 0, 11, 177.9071648067281)
)

$lastRow = 357

foreach ($r in $rows) {
    $rowIndex = $r[0]

    # Carry the date-column formatting (border + centered/top alignment +
    # custom date/time number format) down from the previous row, mirroring
    # what Excel does when a user extends the table downward.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$rowIndex").PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 1).Value = $r[1]
    $ws.Cells.Item($rowIndex, 2).Value = $r[2]
    $ws.Cells.Item($rowIndex, 3).Value = $r[3]
    $ws.Cells.Item($rowIndex, 4).Value = $r[4]

    $lastRow = $rowIndex
}
